# Generate Report for Handoff
# b.md has now been handed off (new handoff file/datetime), updating its
# Status, Latest Handoff File and Latest Handoff Datetime across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) status + handoff date for both locales
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-11 09:27:11"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-11 09:27:08"

# The Latest Handoff File cell (D3) carries a hyperlink whose display text
# must be refreshed to match the new file name. The Hyperlinks collection
# on this host only supports a full-collection rebuild (per-item
# Delete/TextToDisplay edits do not persist), so delete every hyperlink on
# the sheet and recreate them with their original targets, substituting the
# new display text only for D3.
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/a.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd998777571ecccb7a9700490dfdb51e2e27c618/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c0214deb83b8ff9fbd0d1ee06ed10e5943d6f83c/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cac0e040a9826f1bbb708f72f14ccf955e0e6061/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/b.md", "", "", "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/b.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd998777571ecccb7a9700490dfdb51e2e27c618/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c0214deb83b8ff9fbd0d1ee06ed10e5943d6f83c/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cac0e040a9826f1bbb708f72f14ccf955e0e6061/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-11 09:27:11"

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/a.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eeeba86903974fd7da4c4e04e06191ed7a86a379/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/07f1a4724b772b2cf4b49ca6539a5da9e8d92b73/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76c4fd91a247a0c1bdaca3adc379778ec52ca126/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/b.md", "", "", "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f6d1a0186afcb17925ae21e68117abdfe66bc8fe/e2e/b.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eeeba86903974fd7da4c4e04e06191ed7a86a379/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/07f1a4724b772b2cf4b49ca6539a5da9e8d92b73/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76c4fd91a247a0c1bdaca3adc379778ec52ca126/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
